# finished Questions 3 & 4
# Adds VLOOKUP-based lookups under "Question 3" (rows 56-61) and the
# equivalent XLOOKUP-based lookups under "Question 4" (rows 65-70), then
# leaves the sheet scrolled/selected near the Question 4 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
$ws.Activate()

# --- Question 3 : VLOOKUP formulas -----------------------------------
# First row (56) gets its own (non-shared) formula, rows 57:61 are filled
# as one shared-formula group per column, matching how Excel would behave
# if the formula in row 56 were copied down to row 61.
$ws.Range("B56").Formula = '=VLOOKUP(A56,$A$1:$F$52,4)'
$ws.Range("C56").Formula = '=VLOOKUP(A56,$A$1:$P$52,9)'
$ws.Range("D56").Formula = '=VLOOKUP(A56,$A$1:$P$52,14)'

$ws.Range("B57:B61").Formula = '=VLOOKUP(A57,$A$1:$F$52,4)'
$ws.Range("C57:C61").Formula = '=VLOOKUP(A57,$A$1:$P$52,9)'
$ws.Range("D57:D61").Formula = '=VLOOKUP(A57,$A$1:$P$52,14)'

# --- Question 4 : XLOOKUP formulas ------------------------------------
$ws.Range("B65").Formula = '=_xlfn.XLOOKUP(A65,$A$2:$A$52,$D$2:$D$52)'
$ws.Range("C65").Formula = '=_xlfn.XLOOKUP(A65,$A$2:$A$52,$I$2:$I$52)'
$ws.Range("D65").Formula = '=_xlfn.XLOOKUP(A65,$A$2:$A$52,$N$2:$N$52)'

$ws.Range("B66:B70").Formula = '=_xlfn.XLOOKUP(A66,$A$2:$A$52,$D$2:$D$52)'
$ws.Range("C66:C70").Formula = '=_xlfn.XLOOKUP(A66,$A$2:$A$52,$I$2:$I$52)'
$ws.Range("D66:D70").Formula = '=_xlfn.XLOOKUP(A66,$A$2:$A$52,$N$2:$N$52)'

# --- View state : scroll/select near the Question 4 table ------------
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D70").Select()
